# Auto-generated edit script: updates cryptos list data (prices, volumes, and
# re-ordered rows 41-48) per the Fri Mar 8 14:47:39 UTC 2024 GitHub Actions update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '68.340.34'
$ws.Range("E2").Value = '  +1.87%  '

# Row 3
$ws.Range("D3").Value = '3.975.44'
$ws.Range("E3").Value = '  +4.32%  '

# Row 4
$ws.Range("E4").Value = '  -0.27%  '

# Row 5
$ws.Range("D5").Value = '''485.11'
$ws.Range("E5").Value = '  +8.78%  '

# Row 6
$ws.Range("D6").Value = '''149.37'
$ws.Range("E6").Value = '  +2.26%  '

# Row 7
$ws.Range("D7").Value = '''0.628'
$ws.Range("E7").Value = '  +0.91%  '

# Row 8
$ws.Range("E8").Value = '  -0.19%  '

# Row 9
$ws.Range("D9").Value = '''0.735'
$ws.Range("E9").Value = '  -1.30%  '

# Row 10
$ws.Range("D10").Value = '''0.170'
$ws.Range("E10").Value = '  +9.51%  '

# Row 11
$ws.Range("D11").Value = '''0.0000366'
$ws.Range("E11").Value = '  +14.16%  '

# Row 12
$ws.Range("D12").Value = '''43.74'
$ws.Range("E12").Value = '  -0.50%  '

# Row 13
$ws.Range("D13").Value = '4.601.91'
$ws.Range("E13").Value = '  +3.21%  '

# Row 14
$ws.Range("D14").Value = '''10.49'
$ws.Range("E14").Value = '  +0.57%  '

# Row 15
$ws.Range("D15").Value = '''14.91'
$ws.Range("E15").Value = '  +0.94%  '

# Row 16
$ws.Range("D16").Value = '3.974.84'
$ws.Range("E16").Value = '  +3.04%  '

# Row 17
$ws.Range("E17").Value = '  +0.17%  '

# Row 18
$ws.Range("D18").Value = '''19.98'
$ws.Range("E18").Value = '  +0.08%  '

# Row 19
$ws.Range("E19").Value = '  +0.76%  '

# Row 20
$ws.Range("D20").Value = '68.371.65'
$ws.Range("E20").Value = '  +0.99%  '

# Row 21
$ws.Range("D21").Value = '''436.65'
$ws.Range("E21").Value = '  +4.19%  '

# Row 22
$ws.Range("D22").Value = '''3.39'
$ws.Range("E22").Value = '  +4.52%  '

# Row 23
$ws.Range("D23").Value = '''14.51'
$ws.Range("E23").Value = '  -1.03%  '

# Row 24
$ws.Range("D24").Value = '''88.41'
$ws.Range("E24").Value = '  +2.59%  '

# Row 25
$ws.Range("D25").Value = '''3.65'
$ws.Range("E25").Value = '  +6.34%  '

# Row 26
$ws.Range("D26").Value = '''38.83'
$ws.Range("E26").Value = '  +3.98%  '

# Row 27
$ws.Range("D27").Value = '''10.22'
$ws.Range("E27").Value = '  +4.52%  '

# Row 28
$ws.Range("D28").Value = '''9.68'
$ws.Range("E28").Value = '  +3.64%  '

# Row 29
$ws.Range("D29").Value = '''734.03'
$ws.Range("E29").Value = '  +0.02%  '

# Row 30
$ws.Range("D30").Value = '''13.33'
$ws.Range("E30").Value = '  -3.26%  '

# Row 31
$ws.Range("D31").Value = '''0.129'
$ws.Range("E31").Value = '  -2.49%  '

# Row 32
$ws.Range("E32").Value = '  +3.49%  '

# Row 33
$ws.Range("D33").Value = '0.0₃0896'
$ws.Range("E33").Value = '  +32.35%  '

# Row 34
$ws.Range("D34").Value = '''42.06'
$ws.Range("E34").Value = '  -2.86%  '

# Row 35
$ws.Range("D35").Value = '''60.49'
$ws.Range("E35").Value = '  +6.87%  '

# Row 36
$ws.Range("D36").Value = '''0.151'
$ws.Range("E36").Value = '  -4.14%  '

# Row 37
$ws.Range("E37").Value = '  -0.08%  '

# Row 38
$ws.Range("D38").Value = '''5.38'
$ws.Range("E38").Value = '  -2.18%  '

# Row 39
$ws.Range("E39").Value = '  -0.53%  '

# Row 40
$ws.Range("D40").Value = '''3.05'
$ws.Range("E40").Value = '  +4.46%  '

# Row 41
$ws.Range("B41").Value = 'WEMIXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D41").Value = '''2.89'
$ws.Range("E41").Value = '  +8.00%  '

# Row 42
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").Value = '''2.26'
$ws.Range("E42").Value = '  +6.59%  '

# Row 43
$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").Value = '''0.142'
$ws.Range("E43").Value = '  +1.81%  '

# Row 44
$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").Value = '''0.337'
$ws.Range("E44").Value = '  +1.26%  '

# Row 45
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '''1.00'
$ws.Range("E45").Value = '  -0.78%  '

# Row 46
$ws.Range("B46").Value = 'Fetch.AI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D46").Value = '''2.55'
$ws.Range("E46").Value = '  +2.61%  '

# Row 47
$ws.Range("B47").Value = 'LidoDAOToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D47").Value = '''3.43'
$ws.Range("E47").Value = '  +1.85%  '

# Row 48
$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").Value = '''3.26'
$ws.Range("E48").Value = '  +0.68%  '

# Row 49
$ws.Range("D49").Value = '''149.11'
$ws.Range("E49").Value = '  +2.51%  '

# Row 50
$ws.Range("D50").Value = '''2.88'
$ws.Range("E50").Value = '  +0.12%  '

# Row 51
$ws.Range("D51").Value = '''25.27'
$ws.Range("E51").Value = '  -8.52%  '
